# Add two new CFs ("Occupation, lakes, artificial" and "Occupation, water
# courses, artificial") to the CBI (biological footprint) sheet, inserted
# right above the existing "Occupation, shrub land, sclerophyllous" row
# (i.e. as new rows 31 and 32), shifting everything below down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two fresh rows at 31/32 (pushes old row31.. down to row33..) ---
$ws.Rows.Item(31).EntireRow.Insert()
$ws.Rows.Item(31).EntireRow.Insert()

# --- New row 31: Occupation, lakes, artificial ---
$ws.Cells.Item(31, 1).Value = "Occupation, lakes, artificial"
$ws.Cells.Item(31, 2).Value = "natural resource::land"
$ws.Cells.Item(31, 3).Value = [double]"7.6923076923076923E-13"

# --- New row 32: Occupation, water courses, artificial ---
$ws.Cells.Item(32, 1).Value = "Occupation, water courses, artificial"
$ws.Cells.Item(32, 2).Value = "natural resource::land"
$ws.Cells.Item(32, 3).Value = [double]"7.6923076923076923E-13"

# Highlight the two new rows (matches the yellow highlight style already
# used elsewhere in the sheet, e.g. row 14)
$ws.Range("A31:C32").Interior.Color = 65535

# The old highlighted reference row (14) loses its highlight on the first
# two columns but keeps it on the amount column
$ws.Range("A14:B14").Style = "Normal"

# New column D got a width tweak in the original edit (cosmetic, no data)
$ws.Columns.Item(4).ColumnWidth = 11.1666666666667

# Update the active selection / scrolled view to match where the user was
# working (new row 31, column D)
$ws.Range("D31").Select()
